$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.314.20'
$ws.Range('E2').Value = '  -2.68%  '
$ws.Range('D3').Value = '2.997.77'
$ws.Range('E3').Value = '  -3.28%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''584.07'
$ws.Range('E5').Value = '  -1.60%  '
$ws.Range('D6').Value = '''146.30'
$ws.Range('E6').Value = '  -6.45%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -3.22%  '
$ws.Range('D9').Value = '3.001.67'
$ws.Range('E9').Value = '  -3.13%  '
$ws.Range('D10').Value = '''0.148'
$ws.Range('E10').Value = '  -6.31%  '
$ws.Range('D11').Value = '''5.72'
$ws.Range('E11').Value = '  -4.02%  '
$ws.Range('D12').Value = '''0.441'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('E13').Value = '  -4.77%  '
$ws.Range('D14').Value = '''34.58'
$ws.Range('E14').Value = '  -6.55%  '
$ws.Range('D15').Value = '''0.122'
$ws.Range('E15').Value = '  +1.95%  '
$ws.Range('D16').Value = '3.498.99'
$ws.Range('E16').Value = '  -3.04%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '62.419.56'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '''6.99'
$ws.Range('E18').Value = '  -2.97%  '
$ws.Range('D19').Value = '3.004.22'
$ws.Range('E19').Value = '  -3.05%  '
$ws.Range('D20').Value = '''457.94'
$ws.Range('E20').Value = '  -4.69%  '
$ws.Range('D21').Value = '''13.84'
$ws.Range('E21').Value = '  -4.42%  '
$ws.Range('D22').Value = '''0.677'
$ws.Range('E22').Value = '  -5.24%  '
$ws.Range('D23').Value = '''7.34'
$ws.Range('E23').Value = '  -3.01%  '
$ws.Range('D24').Value = '''80.03'
$ws.Range('E24').Value = '  -1.92%  '
$ws.Range('E25').Value = '  -8.51%  '
$ws.Range('D26').Value = '''12.20'
$ws.Range('E26').Value = '  -5.70%  '
$ws.Range('D27').Value = '''10.20'
$ws.Range('E27').Value = '  -4.54%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '''7.16'
$ws.Range('E30').Value = '  -5.94%  '
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('E32').Value = '  -5.71%  '
$ws.Range('D33').Value = '''26.85'
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('E34').Value = '  -5.09%  '
$ws.Range('E35').Value = '  -4.98%  '
$ws.Range('D36').Value = '0.0₃0788'
$ws.Range('E36').Value = '  -6.55%  '
$ws.Range('D37').Value = '''5.72'
$ws.Range('E37').Value = '  -5.33%  '
$ws.Range('E38').Value = '  -6.85%  '
$ws.Range('D39').Value = '''50.01'
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('D40').Value = '''8.90'
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('D41').Value = '''2.93'
$ws.Range('E41').Value = '  -11.06%  '
$ws.Range('D42').Value = '''406.87'
$ws.Range('E42').Value = '  -8.44%  '
$ws.Range('D43').Value = '''0.113'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = '''0.275'
$ws.Range('E44').Value = '  -5.77%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').Value = '''39.36'
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '''0.0352'
$ws.Range('E46').Value = '  -3.66%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.762.80'
$ws.Range('E47').Value = '  -2.65%  '
$ws.Range('D48').Value = '''127.61'
$ws.Range('E48').Value = '  -3.45%  '
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('E50').Value = '  -2.20%  '
$ws.Range('D51').Value = '''23.70'
$ws.Range('E51').Value = '  -9.40%  '
